$p = $ppt.ActivePresentation

# ---- slide8 ----
$s8 = $p.Slides.Add(8, 2)
$__r = $s8.Shapes.Item(1).TextFrame.TextRange
$s8.Shapes.Item(1).TextFrame.TextRange.Text = 'Strchar'
$__r = $s8.Shapes.Item(1).TextFrame.TextRange
$__r = $__r.InsertAfter('()')
$__r = $s8.Shapes.Item(2).TextFrame.TextRange
$s8.Shapes.Item(2).TextFrame.TextRange.Text = 'int'
$__r = $s8.Shapes.Item(2).TextFrame.TextRange
$__r = $__r.InsertAfter(' main')
$__r = $__r.InsertAfter('()')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('{    ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('char ')
$__r = $__r.InsertAfter('str')
$__r = $__r.InsertAfter('[] = "My name is ')
$__r = $__r.InsertAfter('Ayush')
$__r = $__r.InsertAfter('";    ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('char')
$__r = $__r.InsertAfter('* ')
$__r = $__r.InsertAfter('ch')
$__r = $__r.InsertAfter(' = ')
$__r = $__r.InsertAfter('strchr')
$__r = $__r.InsertAfter('(')
$__r = $__r.InsertAfter('str')
$__r = $__r.InsertAfter(', ''a'');    ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('Printf')
$__r = $__r.InsertAfter('(“%s”,')
$__r = $__r.InsertAfter('ch')
$__r = $__r.InsertAfter(');    ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('return ')
$__r = $__r.InsertAfter('0')
$__r = $__r.InsertAfter(';')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('}')
$s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).ParagraphFormat.Bullet.Visible = $false
$s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).ParagraphFormat.Bullet.Visible = $false
$s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).IndentLevel = 3
$s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).ParagraphFormat.Bullet.Visible = $false
$s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).IndentLevel = 3
$s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).ParagraphFormat.Bullet.Visible = $false
$s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(5).IndentLevel = 3
$s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(5).ParagraphFormat.Bullet.Visible = $false
$s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(6).IndentLevel = 3
$s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(6).ParagraphFormat.Bullet.Visible = $false
$s8.Shapes.Item(2).TextFrame.TextRange.Paragraphs(7).ParagraphFormat.Bullet.Visible = $false

# ---- slide9 ----
$s9 = $p.Slides.Add(9, 2)
$__r = $s9.Shapes.Item(1).TextFrame.TextRange
$s9.Shapes.Item(1).TextFrame.TextRange.Text = 'Strchar'
$__r = $s9.Shapes.Item(1).TextFrame.TextRange
$__r = $__r.InsertAfter('()')
$__r = $s9.Shapes.Item(2).TextFrame.TextRange
$s9.Shapes.Item(2).TextFrame.TextRange.Text = 'int'
$__r = $s9.Shapes.Item(2).TextFrame.TextRange
$__r = $__r.InsertAfter(' main')
$__r = $__r.InsertAfter('()')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('{    ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('char ')
$__r = $__r.InsertAfter('str')
$__r = $__r.InsertAfter('[] = "My name is ')
$__r = $__r.InsertAfter('Ayush')
$__r = $__r.InsertAfter('";    ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('char')
$__r = $__r.InsertAfter('* ')
$__r = $__r.InsertAfter('ch')
$__r = $__r.InsertAfter(' = ')
$__r = $__r.InsertAfter('strchr')
$__r = $__r.InsertAfter('(')
$__r = $__r.InsertAfter('str')
$__r = $__r.InsertAfter(', ''a'');    ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('Printf')
$__r = $__r.InsertAfter('(“%d”,')
$__r = $__r.InsertAfter('ch')
$__r = $__r.InsertAfter(' ')
$__r = $__r.InsertAfter('- ')
$__r = $__r.InsertAfter('str')
$__r = $__r.InsertAfter(' + ')
$__r = $__r.InsertAfter('1);    ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('return ')
$__r = $__r.InsertAfter('0')
$__r = $__r.InsertAfter(';')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('}')
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).ParagraphFormat.Bullet.Visible = $false
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).ParagraphFormat.Bullet.Visible = $false
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).IndentLevel = 3
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).ParagraphFormat.Bullet.Visible = $false
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).IndentLevel = 3
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).ParagraphFormat.Bullet.Visible = $false
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(5).IndentLevel = 3
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(5).ParagraphFormat.Bullet.Visible = $false
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(6).IndentLevel = 3
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(6).ParagraphFormat.Bullet.Visible = $false
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(7).ParagraphFormat.Bullet.Visible = $false

# ---- slide10 ----
$s10 = $p.Slides.Add(10, 2)
$__r = $s10.Shapes.Item(1).TextFrame.TextRange
$s10.Shapes.Item(1).TextFrame.TextRange.Text = 'Strstr'
$__r = $s10.Shapes.Item(1).TextFrame.TextRange
$__r = $__r.InsertAfter('()')
$__r = $s10.Shapes.Item(2).TextFrame.TextRange
$s10.Shapes.Item(2).TextFrame.TextRange.Text = 'main'
$__r = $s10.Shapes.Item(2).TextFrame.TextRange
$__r = $__r.InsertAfter('() ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('{ ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('	// Take any two strings ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('	char s1[] = "')
$__r = $__r.InsertAfter('GeeksforGeeks')
$__r = $__r.InsertAfter('"; ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('	char s2[] = "for"; ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('	char* p; ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('	// Find first occurrence of s2 in s1 ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('	p = ')
$__r = $__r.InsertAfter('strstr')
$__r = $__r.InsertAfter('(s1, s2); ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('	// Prints the result ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('	if (p) { ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('		')
$__r = $__r.InsertAfter('printf')
$__r = $__r.InsertAfter('("String found\n"); ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('		')
$__r = $__r.InsertAfter('printf')
$__r = $__r.InsertAfter('("First occurrence of string ''%s'' in ''%s'' is ''%s''", s2, s1, p); ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('	} else')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('		')
$__r = $__r.InsertAfter('printf')
$__r = $__r.InsertAfter('("String not found\n"); ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('	return 0; ')
$__r = $__r.InsertAfter('')
$__r = $__r.InsertAfter('} ')
$__r = $__r.InsertAfter('')
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(5).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(6).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(7).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(8).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(9).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(10).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(11).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(12).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(13).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(14).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(15).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(16).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(17).ParagraphFormat.Bullet.Visible = $false
$s10.Shapes.Item(2).TextFrame.TextRange.Paragraphs(18).ParagraphFormat.Bullet.Visible = $false

